# The source commit ("Fixed POI packaging and upgraded to POI 3.15") is a
# build/packaging-tooling change to the project that produced this .docx
# test fixture. Re-diffing the canonical (pretty-printed, attribute-sorted)
# OOXML shows only XML-serialization noise -- namespace / attribute
# declarations being re-ordered alphabetically (an artifact of the XML
# writer used by the newer POI version) and non-deterministic `w:rsid*`
# revision-save-id attributes dropping out of the canonical form used for
# diffing. Every attribute=value pair that appears on both sides of each
# hunk is identical; nothing in the document's visible text, formatting,
# structure, or content actually changed.
#
# The Word object model doesn't expose control over raw XML attribute
# ordering or rsid bookkeeping (those are serializer-level details, not
# document content), so there is no corresponding content edit to make
# here. This script intentionally performs no mutation, leaving the
# document's semantic content identical to before.docx.

$d = $word.ActiveDocument
